$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Outbound clicks" header to "Unique outbound clicks"
$ws.Range("M1").Value = "Unique outbound clicks"

# The renamed cell picks up its own font (Calibri 11, black) ...
$ws.Range("M1").Font.Name = "Calibri"
$ws.Range("M1").Font.Size = 11
$ws.Range("M1").Font.Color = 0

# ... and the header row grows slightly taller to fit the new text
$ws.Rows(1).RowHeight = 36.5

# Move the active selection to N5
$ws.Range("N5").Select()
